$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (3-6) get rotated: the former row 6 data moves up to
# row 3, and rows 3, 4, 5 each shift down by one row. Only columns D and
# L:T (Fecha, Calidad .. Kg/unidad) change; A:K stay the same per row.

# --- Row 3 <= old Row 6 ---
$ws.Range("D3").Value = 44334
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 12500
$ws.Range("Q3").Value = "$/caja 12 kilos empedrada"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1042
$ws.Range("T3").Value = 12

# --- Row 4 <= old Row 3 ---
$ws.Range("D4").Value = 44316
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 17500
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17750
$ws.Range("Q4").Value = "$/caja 16 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1109
$ws.Range("T4").Value = 16

# --- Row 5 <= old Row 4 ---
$ws.Range("D5").Value = 44316
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("Q5").Value = "$/caja 16 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 16

# --- Row 6 <= old Row 5 ---
$ws.Range("D6").Value = 44330
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18
